# Add a new bug-tracker row (issue #63549 - tensorflow/XLA) to the "Tensorflow" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tensorflow")

# --- Values -------------------------------------------------------------
$ws.Range("A2").Value = 63549
$ws.Range("B2").Value = 7899
$ws.Range("C2").Value = "Inconsistency in XLA Cotionmpila with Operand Order Swap in tf.add with Specific Operators on GPU"
$ws.Range("D2").Value = "Yes"
$ws.Range("E2").Value = "Yes"
$ws.Range("F2").Value = "https://github.com/openxla/xla/tree/main/xla/service/gpu/ir_emission_utils.cc"
$ws.Range("G2").Value = "FindTransposeHero, transpose"

# --- Hyperlinks (Issue #, PR # and Buggy File link to their URLs) -------
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/tensorflow/tensorflow/issues/63549")
$ws.Hyperlinks.Add($ws.Range("B2"), "https://github.com/openxla/xla/pull/7899")
$ws.Hyperlinks.Add($ws.Range("F2"), "https://github.com/openxla/xla/tree/main/xla/service/gpu/ir_emission_utils.cc")

# --- Re-apply the sheet's existing cell styles ---------------------------
# Hyperlinks.Add() stamps its own built-in "Hyperlink" style; restore the
# workbook's existing hyperlink look (underlined blue, same as the other
# sheets' A/B columns) by pasting formats from a cell that already carries
# that exact style.
$wb.Worksheets.Item("Jax").Range("A2").Copy()
$ws.Range("A2").PasteSpecial(-4122)
$ws.Range("B2").PasteSpecial(-4122)
$ws.Range("F2").PasteSpecial(-4122)

# Plain text cells (Title / Reproduced / Fixed / Buggy Function(s)) share
# the same style as the header row - copy it from there.
$ws.Range("C1").Copy()
$ws.Range("C2").PasteSpecial(-4122)
$ws.Range("D2").PasteSpecial(-4122)
$ws.Range("E2").PasteSpecial(-4122)
$ws.Range("G2").PasteSpecial(-4122)

# --- Column width for the new Title column -------------------------------
$ws.Columns.Item(3).ColumnWidth = 75.5

Write-Output "Added bug from issue 63549 - tensorflow"
